$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Extend formatting from column I into the new columns J:N (years 2019-2023) ---
# Row 5 has no D:I cells at all (plain section caption), so it is skipped -
# only rows 3, 4 and 6-19 in column I carry the borders/number-format that the
# new year columns must inherit. Rows 20-29 stay untouched (no new cells).
$ws.Range("I3:I4").Copy()
$ws.Range("J3:N4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I6:I19").Copy()
$ws.Range("J6:N19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 2. Year header row (row 4) ---
$ws.Cells.Item(4, 10).Value = 2019
$ws.Cells.Item(4, 11).Value = 2020
$ws.Cells.Item(4, 12).Value = 2021
$ws.Cells.Item(4, 13).Value = 2022
$ws.Cells.Item(4, 14).Value = 2023

# --- 3. Data values for the new years, one row at a time ---
$rowData = @{
  6  = @(81.1, 85.8, 78.1, 72.2, 75.7)
  7  = @(18.9, 14.2, 21.9, 27.8, 24.3)
  9  = @(22.8, 25.6, 24.2, 21.4, 31.1)
  10 = @(77.2, 74.4, 75.8, 78.6, 68.9)
  12 = @(84.4, 72.7, 73.3, 72.8, 76.7)
  13 = @(15.6, 27.3, 26.7, 27.2, 23.3)
  15 = @(90.3, 93.4, 90.5, 87.8, 89)
  16 = @(9.7, 6.6, 9.5, 12.2, 11)
  18 = @(60.2, 66, 59.3, 44.9, 48.3)
  19 = @(39.8, 34, 40.7, 55.1, 51.7)
}

foreach ($r in $rowData.Keys) {
  $vals = $rowData[$r]
  for ($i = 0; $i -lt $vals.Length; $i++) {
    $col = 10 + $i   # column J = 10 ... N = 14
    $ws.Cells.Item($r, $col).Value = $vals[$i]
  }
}

# Rows 8, 11, 14, 17 (sub-header rows) already received the blank, bordered
# style from the paste above, so they need no values - matches the diff.

# --- 4. Row 20 restyle: new 8pt Times New Roman font + taller row ---
$ws.Range("A20:C20").Font.Size = 8
$ws.Rows.Item(20).RowHeight = 15.75
